$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new rows 41 and 42 (practice room / play offline strings) ---
# Carry over the same cell formatting (border + wrap + vertical-center)
# already used by the existing data rows.
$ws.Range("A40:G40").Copy()
$ws.Range("A41:G42").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Cell values are written in the same order the shared-string table was
# originally built in, so new shared-string indices line up with the
# reference workbook.

# Row 41: UI Scene name / tmp practice room / Practice Room / 練習室 / 练习室 / Salle de Test / Sala de Prueba
$ws.Cells.Item(41, 1).Value = "UI Scene name"
$ws.Cells.Item(41, 2).Value = "tmp practice room"
$ws.Cells.Item(41, 3).Value = "Practice Room"
$ws.Cells.Item(41, 6).Value = "練習室"
$ws.Cells.Item(41, 7).Value = "练习室"
$ws.Cells.Item(41, 4).Value = "Salle de Test"
$ws.Cells.Item(41, 5).Value = "Sala de Prueba"

# Row 42: UI button / btn play offline / Play offline / オフライン / Juego local / Jouer localement / 本地播放
$ws.Cells.Item(42, 1).Value = "UI button"
$ws.Cells.Item(42, 2).Value = "btn play offline"
$ws.Cells.Item(42, 3).Value = "Play offline"
$ws.Cells.Item(42, 6).Value = "オフライン"
$ws.Cells.Item(42, 5).Value = "Juego local"
$ws.Cells.Item(42, 4).Value = "Jouer localement"
$ws.Cells.Item(42, 7).Value = "本地播放"

$ws.Rows.Item(42).RowHeight = 28.8

# --- View / window state ---
# Scroll the sheet so row 30 is at the top, then select F44 (the cell the
# author was last looking at / editing next to the new rows).
$ws.Select()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F44").Select()
